# Update row 5 values to the new "custom accuracy" readings, then remove
# row 6 entirely (data trimmed from 2 rows to 1 row of readings).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row5 = @{
    "B5"  = 13.11
    "C5"  = 9.949999999999999
    "D5"  = 0.61
    "E5"  = 28.05
    "F5"  = 23.58
    "G5"  = 10.07
    "H5"  = 42.68
    "I5"  = 15.45
    "J5"  = 7.05
    "K5"  = 10.71
    "L5"  = 11.4
    "M5"  = 11.91
    "N5"  = 3.24
    "O5"  = 10.03
    "P5"  = 14.5
    "Q5"  = 8.380000000000001
    "R5"  = 0.42
    "S5"  = 0.26
    "T5"  = 147.02
    "U5"  = 28.27
    "V5"  = 9.220000000000001
    "W5"  = 19.11
    "X5"  = 10.28
    "Y5"  = 1.34
    "Z5"  = 20.98
    "AA5" = 7.96
    "AB5" = 7.77
    "AC5" = 8.369999999999999
    "AD5" = 11.99
    "AE5" = 0
    "AF5" = 38.61
    "AG5" = 5.49
    "AH5" = 11.55
}

foreach ($addr in $row5.Keys) {
    $ws.Range($addr).Value = $row5[$addr]
}

# Drop the now-unused sixth data row; the sheet's used range becomes
# A1:AH5 afterwards.
$ws.Rows.Item(6).Delete()
